# B6-PowerPoint.pptx edit: Sat, Jun 13, 2020 10:04:40 AM
#
# 1) Three tables (on slides 14, 15, 16) switch from the custom
#    "Table_0" table style ({44311CAA-88B3-4980-9638-23AC147A936B}) to the
#    built-in table style {6138919E-21F6-4239-A4B3-49B48093B712}.
# 2) The deck's working theme (applied to the slide master, and so to every
#    slide) has its 12 theme colors changed from the "Red Violet"/Integral
#    palette back to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables ---------------------------------------
$newTableStyle = "{6138919E-21F6-4239-A4B3-49B48093B712}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Swap the theme color palette ------------------------------------
# Office (stock) theme colors, in clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = 0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
